# Weekly data refresh: a new sampling row is prepended to the price table
# (row 31), pushing every subsequent row down by one. The oldest row that
# falls off the bottom (old row 146) re-appears as the new last row (147).
#
# Inserting a whole row at 31 shifts rows 31:146 down to 32:147 in one shot
# (carrying along each cell's value *and* formatting, e.g. the custom date
# style on column D), which reproduces that shift exactly. All that is left
# is to populate the freshly inserted (empty) row 31 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 4
$ws.Range("B31").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C31").Value = "Los Lagos"
$ws.Range("D31").Value = 44608
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = 100112009
$ws.Range("G31").Value = "Acelga"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = 10000
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = 10000
$ws.Range("N31").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O31").Value = "Región de La Araucanía"
$ws.Range("P31").Value = 833
$ws.Range("Q31").Value = 12
$ws.Range("R31").Value = "Hortaliza"
